$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add "Instructor titular:" label in D5 (bold, like the other labels in column A)
#    and the instructor name in E5.
$ws.Range("D5").Value = "Instructor titular:"
$ws.Range("D5").Font.Bold = $true
$ws.Range("E5").Value = "LILIANA MARÍA GALEANO ZEA "

# 2. Fix spacing in "Proyecto 2 + Manual Técnico" -> "Proyecto 2  + Manual Técnico"
#    (double space after "2") for the schedule cells referencing
#    YEISON BARRIOS FUNIELES, and for the keyword table entry in row 28.
$proyecto2New = "Proyecto 2  + Manual Técnico `n YEISON BARRIOS FUNIELES `n 801"

$cells = @("B9", "C9", "E9", "F9", "B10", "C10", "E10", "F10")
foreach ($cellAddr in $cells) {
    $ws.Range($cellAddr).Value = $proyecto2New
}

# Editing the wrapped text above can trigger an automatic row-height
# resize; restore the rows to their original auto-fit (no custom height).
$ws.Rows(9).EntireRow.AutoFit()
$ws.Rows(10).EntireRow.AutoFit()

$ws.Range("A28").Value = "Proyecto 2  + Manual Técnico"
